$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price and 1h volume/change figures
$ws.Range('D2').Value = '56.520.77'
$ws.Range('E2').Value = '  +4.06%  '
$ws.Range('D3').Value = '2.995.06'
$ws.Range('E3').Value = '  +4.24%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '507.22'
$ws.Range('E5').Value = '  +6.81%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '136.59'
$ws.Range('E6').Value = '  +8.32%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.433'
$ws.Range('E8').Value = '  +7.51%  '
$ws.Range('E9').Value = '  +12.57%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.109'
$ws.Range('E10').Value = '  +12.89%  '
$ws.Range('E11').Value = '  +7.16%  '
$ws.Range('E12').Value = '  +4.65%  '
$ws.Range('D13').Value = '3.502.00'
$ws.Range('E13').Value = '  +4.57%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.91'
$ws.Range('E15').Value = '  +14.84%  '
$ws.Range('D16').Value = '56.530.93'
$ws.Range('E16').Value = '  +4.28%  '
$ws.Range('D17').Value = '2.987.04'
$ws.Range('E17').Value = '  +4.39%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.82'
$ws.Range('E18').Value = '  +9.95%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.49'
$ws.Range('E19').Value = '  +8.73%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.82'
$ws.Range('E20').Value = '  +10.44%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '326.85'
$ws.Range('E21').Value = '  +10.84%  '
$ws.Range('E22').Value = '  +0.03%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.476'
$ws.Range('E23').Value = '  +7.15%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '62.33'
$ws.Range('E24').Value = '  +6.05%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.00'
$ws.Range('E25').Value = '  -0.17%  '
$ws.Range('E26').Value = '  +6.63%  '
$ws.Range('D27').Value = '0.0₃0907'
$ws.Range('E27').Value = '  +11.95%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.49'
$ws.Range('E28').Value = '  +3.29%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.94'
$ws.Range('E29').Value = '  +12.23%  '
$ws.Range('E30').Value = '  +7.74%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.77'
$ws.Range('E31').Value = '  +9.37%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.60'
$ws.Range('E32').Value = '  +8.42%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '157.26'
$ws.Range('E33').Value = '  +16.71%  '
$ws.Range('E34').Value = '  +6.14%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.59'
$ws.Range('E35').Value = '  +3.31%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.27'
$ws.Range('E36').Value = '  +4.76%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0677'
$ws.Range('E37').Value = '  +10.10%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '23.47'
$ws.Range('E38').Value = '  +2.35%  '
$ws.Range('D39').Value = '3.025.23'
$ws.Range('E39').Value = '  +4.67%  '
$ws.Range('E40').Value = '  +3.12%  '
$ws.Range('E41').Value = '  +0.36%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.644'
$ws.Range('E42').Value = '  +7.49%  '
$ws.Range('D43').Value = '2.260.69'
$ws.Range('E43').Value = '  +10.69%  '
$ws.Range('E44').Value = '  +7.43%  '
$ws.Range('E45').Value = '  +2.74%  '
$ws.Range('E46').Value = '  +5.98%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.97'
$ws.Range('E47').Value = '  +22.12%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0237'
$ws.Range('E48').Value = '  +11.95%  '
$ws.Range('E49').Value = '  +8.76%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '19.14'
$ws.Range('E50').Value = '  +7.16%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0871'
$ws.Range('E51').Value = '  +10.99%  '
